$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.838.53"
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("D3").Value = "3.430.76"
$ws.Range("E3").Value = "  +4.69%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'256.31"
$ws.Range("E5").Value = "  +1.62%  "
$ws.Range("D6").Value = "'658.59"
$ws.Range("E6").Value = "  +6.06%  "
$ws.Range("E7").Value = "  +5.09%  "
$ws.Range("D8").Value = "'0.432"
$ws.Range("E8").Value = "  +7.95%  "
$ws.Range("E9").Value = "  +10.72%  "
$ws.Range("D10").Value = "'0.999"
$ws.Range("E10").Value = "  -0.04%  "
$ws.Range("D11").Value = "3.428.21"
$ws.Range("E11").Value = "  +4.75%  "
$ws.Range("E12").Value = "  +7.21%  "
$ws.Range("D13").Value = "'42.03"
$ws.Range("E13").Value = "  +8.01%  "
$ws.Range("D14").Value = "'6.41"
$ws.Range("E14").Value = "  +17.60%  "
$ws.Range("E15").Value = "  +5.40%  "
$ws.Range("D16").Value = "97.479.26"
$ws.Range("E16").Value = "  -1.33%  "
$ws.Range("D17").Value = "4.066.88"
$ws.Range("E17").Value = "  +4.28%  "
$ws.Range("D18").Value = "'8.74"
$ws.Range("E18").Value = "  +39.21%  "
$ws.Range("D19").Value = "3.444.40"
$ws.Range("E19").Value = "  +5.11%  "
$ws.Range("D20").Value = "'17.61"
$ws.Range("E20").Value = "  +15.19%  "
$ws.Range("D21").Value = "'0.509"
$ws.Range("E21").Value = "  +61.59%  "
$ws.Range("E22").Value = "  +20.18%  "
$ws.Range("D23").Value = "'3.47"
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("D24").Value = "'509.84"
$ws.Range("E24").Value = "  +5.43%  "
$ws.Range("E25").Value = "  +4.05%  "
$ws.Range("D26").Value = "'6.20"
$ws.Range("E26").Value = "  +10.74%  "
$ws.Range("D27").Value = "'98.99"
$ws.Range("E27").Value = "  +11.52%  "
$ws.Range("D28").Value = "'12.88"
$ws.Range("E28").Value = "  +8.19%  "
$ws.Range("D29").Value = "'0.155"
$ws.Range("E29").Value = "  +12.20%  "
$ws.Range("D30").Value = "'11.47"
$ws.Range("E30").Value = "  +12.16%  "
$ws.Range("D31").Value = "'0.200"
$ws.Range("E31").Value = "  +5.86%  "
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("D34").Value = "'0.575"
$ws.Range("E34").Value = "  +22.24%  "
$ws.Range("D35").Value = "'30.16"
$ws.Range("E35").Value = "  +8.62%  "
$ws.Range("E36").Value = "  +17.56%  "
$ws.Range("D37").Value = "'7.87"
$ws.Range("E37").Value = "  +9.94%  "
$ws.Range("E38").Value = "  +5.37%  "
$ws.Range("D39").Value = "'1.43"
$ws.Range("E39").Value = "  +16.75%  "
$ws.Range("D40").Value = "'518.57"
$ws.Range("E40").Value = "  +6.37%  "
$ws.Range("D41").Value = "'24.74"
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("D42").Value = "'0.872"
$ws.Range("E42").Value = "  +13.56%  "
$ws.Range("E43").Value = "  +2.64%  "
$ws.Range("D44").Value = "'0.0415"
$ws.Range("E44").Value = "  +24.16%  "
$ws.Range("D45").Value = "'3.32"
$ws.Range("E45").Value = "  +8.20%  "
$ws.Range("D46").Value = "'5.49"
$ws.Range("E46").Value = "  +17.54%  "
$ws.Range("B47").Value = "Cosmos"
$ws.Range("C47").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D47").Value = "'8.21"
$ws.Range("E47").Value = "  +14.01%  "
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").Value = "'1.00"
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("E49").Value = "  +17.58%  "
$ws.Range("E50").Value = "  +8.15%  "
$ws.Range("D51").Value = "'51.45"
$ws.Range("E51").Value = "  +11.10%  "
